# SOURCE_C.xlsx fixture fix: "id" -> "ID" header, and active-sheet/selection
# cleanup (C.1 becomes the active tab with A1 selected; C.PIVOT's selection
# resets to A1 as well).

$wb = $excel.ActiveWorkbook

$wsC1 = $wb.Worksheets.Item("C.1")
$wsPivot = $wb.Worksheets.Item("C.PIVOT")

# Fix the header text casing on the C.1 sheet.
$wsC1.Range("A1").Value = "ID"

# C.PIVOT's remembered selection resets to A1 (was C1).
$wsPivot.Activate()
[void]$wsPivot.Range("A1").Select()

# Make C.1 the active sheet/tab, with A1 selected (was A2).
$wsC1.Activate()
[void]$wsC1.Range("A1").Select()
